$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the existing data (rows 1-2 shift down to rows 5-6)
$ws.Rows("1:4").Insert()

# Fill in the new category rows, matching the shared-string insertion order
# observed in the target workbook: A1, B1, A2, A3, A4, B2, B3, B4
$ws.Range("A1").Value = "Beginner"
$ws.Range("B1").Value = "Tree Pose 1,Tree Pose 2"

$ws.Range("A2").Value = "Intermediate"
$ws.Range("A3").Value = "Advanced"
$ws.Range("A4").Value = "Expert"

$ws.Range("B2").Value = "Tree Pose 3,Tree Pose 4"
$ws.Range("B3").Value = "Tree Pose 5,Tree Pose 6"
$ws.Range("B4").Value = "Tree Pose 7,Tree Pose 8"

# The hyperlink on the video-path cell should follow it down to its new row:
# remove the stale link (still anchored at the pre-insert B1 position) and
# re-create it on B5, then restore the Hyperlink cell style (Add() mints a
# fresh style entry, so point the cell back at the original named style).
$ws.Range("B5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B5"), "file:///\\Videos\TreePose.mp4") | Out-Null
$ws.Range("B5").Style = "Hyperlink"

# Update selection to match the target state
$ws.Range("M6").Select()
